$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatted "A" column style (bold/centered/bordered header style)
# down into the two brand-new rows before writing their values, so the
# cellXfs table isn't perturbed with extra intermediate entries.
$ws.Range("A15").Copy($ws.Range("A16:A17"))

# Rows 8-17, columns A-E (A=id, B=name, C=from_bus, D=to_bus, E=in_service).
# Rows 8 & 9 are renamed from extr1/extr2 -> line7/line8 (values also change).
# Rows 10-15 are renamed extr3..extr8 -> extr1..extr6 (shifted up), values change.
# Rows 16 & 17 are brand new: extr7, extr8.
$data = @(
    @(8,  6,  "line7", 14, 11, $false),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $true),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $false),
    @(15, 13, "extr6", 7,  11, $true),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $data) {
    $r = $row[0]
    $idVal = $row[1]
    $name = $row[2]
    $cVal = $row[3]
    $dVal = $row[4]
    $eVal = $row[5]

    $ws.Cells.Item($r, 1).Value = $idVal
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $cVal
    $ws.Cells.Item($r, 4).Value = $dVal
    $ws.Cells.Item($r, 5).Value = $eVal
}
